$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Next empty row after the existing data (row 26 -> new row 27)
$row = $ws.UsedRange.Rows.Count + 1

# Text-like columns (Date / Time / Weekday / Week). Date ("2025-01-15") and
# Week ("02") look numeric to Excel's auto-detection, so force them in as
# literal text via a leading apostrophe, then strip the resulting
# "quote prefix" style back off so the cell keeps the workbook's default
# (unstyled) formatting - matching every other data row.
$ws.Cells.Item($row, 1).Value = "'2025-01-15"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "20:00:52"

$ws.Cells.Item($row, 3).Value = "Wednesday"

$ws.Cells.Item($row, 4).Value = "'02"
$ws.Cells.Item($row, 4).Style = "Normal"

# Numeric resale-count columns (Beijing .. Wuhan)
$ws.Cells.Item($row, 5).Value = 126882
$ws.Cells.Item($row, 6).Value = 141303
$ws.Cells.Item($row, 7).Value = 169408
$ws.Cells.Item($row, 8).Value = 157522
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142933
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191325
$ws.Cells.Item($row, 14).Value = 115439
$ws.Cells.Item($row, 15).Value = 45126
$ws.Cells.Item($row, 16).Value = 28525
$ws.Cells.Item($row, 17).Value = 65701
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 49291
$ws.Cells.Item($row, 20).Value = -1
